# Implement counting of user response.
# Update existing chat-history rows' timestamps/values and append three new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 1-3 ---
$ws.Range("B1").Value = "18:01:16"
$ws.Range("D1").Value = 741216406

$ws.Range("B2").Value = "18:01:18"

$ws.Range("B3").Value = "18:01:27"
$ws.Range("C3").Value = -485430438

# --- Append new rows 4-6 ---
# Force text format on column A so the date-looking string "2021-11-06"
# is kept as literal text (matching the other rows) instead of being
# auto-converted into a date serial number, then clear the formatting
# again so no extra cell style lingers behind.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2021-11-06"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = "18:01:30"
$ws.Range("C4").Value = -615761128
$ws.Range("D4").Value = 741216406

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2021-11-06"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "18:01:40"
$ws.Range("C5").Value = -615761128
$ws.Range("D5").Value = 1107423707

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2021-11-06"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = "18:01:42"
$ws.Range("C6").Value = -485430438
$ws.Range("D6").Value = 1107423707
